$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.295.72"
$ws.Cells.Item(2, 5).Value = "  -1.53%  "
$ws.Cells.Item(3, 4).Value = "2.281.35"
$ws.Cells.Item(3, 5).Value = "  -1.30%  "
$ws.Cells.Item(4, 5).Value = "  -0.20%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "112.41"
$ws.Cells.Item(5, 5).Value = "  -2.30%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "265.55"
$ws.Cells.Item(6, 5).Value = "  -1.43%  "
$ws.Cells.Item(7, 5).Value = "  -0.86%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.607"
$ws.Cells.Item(9, 5).Value = "  -3.16%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "47.70"
$ws.Cells.Item(10, 5).Value = "  -2.22%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0932"
$ws.Cells.Item(11, 5).Value = "  -1.49%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "8.85"
$ws.Cells.Item(12, 5).Value = "  +0.69%  "
$ws.Cells.Item(13, 5).Value = "  +0.96%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "15.53"
$ws.Cells.Item(14, 5).Value = "  -1.08%  "
$ws.Cells.Item(15, 4).Value = "2.622.91"
$ws.Cells.Item(15, 5).Value = "  -0.10%  "
$ws.Cells.Item(16, 5).Value = "  -0.75%  "
$ws.Cells.Item(17, 4).Value = "2.276.65"
$ws.Cells.Item(17, 5).Value = "  -1.75%  "
$ws.Cells.Item(18, 4).Value = "43.234.71"
$ws.Cells.Item(18, 5).Value = "  -1.31%  "
$ws.Cells.Item(19, 5).Value = "  -1.94%  "
$ws.Cells.Item(20, 5).Value = "  +2.75%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "71.50"
$ws.Cells.Item(21, 5).Value = "  -1.54%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2.51"
$ws.Cells.Item(22, 5).Value = "  -2.21%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "232.15"
$ws.Cells.Item(23, 5).Value = "  -1.02%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "9.64"
$ws.Cells.Item(24, 5).Value = "  +1.29%  "
$ws.Cells.Item(25, 5).Value = "  -0.23%  "
$ws.Cells.Item(26, 5).Value = "  +0.80%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.33"
$ws.Cells.Item(27, 5).Value = "  -2.00%  "
$ws.Cells.Item(28, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "40.48"
$ws.Cells.Item(28, 5).Value = "  -8.25%  "
$ws.Cells.Item(29, 2).Value = "WEMIXToken"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "3.35"
$ws.Cells.Item(29, 5).Value = "  -1.96%  "
$ws.Cells.Item(30, 2).Value = "Toncoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.25"
$ws.Cells.Item(30, 5).Value = "  -1.17%  "
$ws.Cells.Item(31, 2).Value = "Monero"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "172.08"
$ws.Cells.Item(31, 5).Value = "  -3.28%  "
$ws.Cells.Item(32, 2).Value = "EthereumClassic"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "21.35"
$ws.Cells.Item(32, 5).Value = "  -2.56%  "
$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0908"
$ws.Cells.Item(33, 5).Value = "  -3.09%  "
$ws.Cells.Item(34, 2).Value = "Filecoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.82"
$ws.Cells.Item(34, 5).Value = "  +4.58%  "
$ws.Cells.Item(35, 2).Value = "Stellar"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.127"
$ws.Cells.Item(35, 5).Value = "  +0.41%  "
$ws.Cells.Item(36, 2).Value = "RenderToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.64"
$ws.Cells.Item(36, 5).Value = "  -2.35%  "
$ws.Cells.Item(37, 2).Value = "NEARProtocol"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.92"
$ws.Cells.Item(37, 5).Value = "  -0.86%  "
$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0355"
$ws.Cells.Item(38, 5).Value = "  -0.29%  "
$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.104"
$ws.Cells.Item(39, 5).Value = "  -5.53%  "
$ws.Cells.Item(40, 2).Value = "LidoDAOToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.66"
$ws.Cells.Item(40, 5).Value = "  +11.01%  "
$ws.Cells.Item(41, 2).Value = "MultiversX"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "77.06"
$ws.Cells.Item(41, 5).Value = "  +1.28%  "
$ws.Cells.Item(42, 2).Value = "Celestia"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "13.85"
$ws.Cells.Item(42, 5).Value = "  +4.79%  "
$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.237"
$ws.Cells.Item(43, 5).Value = "  -4.51%  "
$ws.Cells.Item(44, 2).Value = "THORChain"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "6.26"
$ws.Cells.Item(44, 5).Value = "  +3.62%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.00"
$ws.Cells.Item(45, 5).Value = "  +0.03%  "
$ws.Cells.Item(46, 2).Value = "ARBITRUM"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.39"
$ws.Cells.Item(46, 5).Value = "  -0.19%  "
$ws.Cells.Item(47, 2).Value = "FraxShare"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "8.68"
$ws.Cells.Item(47, 5).Value = "  -1.51%  "
$ws.Cells.Item(48, 2).Value = "Aave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "104.16"
$ws.Cells.Item(48, 5).Value = "  +2.97%  "
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0994"
$ws.Cells.Item(49, 5).Value = "  -1.85%  "
$ws.Cells.Item(50, 2).Value = "TrustWalletToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.25"
$ws.Cells.Item(50, 5).Value = "  +1.21%  "
$ws.Cells.Item(51, 2).Value = "WOONetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.441"
$ws.Cells.Item(51, 5).Value = "  -4.58%  "
